$d = $word.ActiveDocument

# 1. Merge runs around "between different school shootings" (no visible text change,
#    just collapses the run that split "different " from its neighbours).
$d.Content.Find.Execute(
  "There was no way to draw edges between different school shootings, resulting in a lack of directionality.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "There was no way to draw edges between different school shootings, resulting in a lack of directionality.",
  2) | Out-Null

# 2. Merge "...14th Parliament." + " " into a single run with trailing space.
$d.Content.Find.Execute(
  "asked in Singapore’s 14th Parliament. Finding this document",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "asked in Singapore’s 14th Parliament. Finding this document",
  2) | Out-Null

# 3. Merge "...access past " + "OPs" into one run.
$d.Content.Find.Execute(
  "This reality meant I had to be more innovative to access past OPs.",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "This reality meant I had to be more innovative to access past OPs.",
  2) | Out-Null

# 4. Merge " initially worked for a couple of OPs. " + "However, I soon found out that Parliament S"
$d.Content.Find.Execute(
  "initially worked for a couple of OPs. However, I soon found out that Parliament Sessions",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "initially worked for a couple of OPs. However, I soon found out that Parliament Sessions",
  2) | Out-Null

# 5. Merge "...error process. " + "The most optimal method..."
$d.Content.Find.Execute(
  "a time-consuming trial and error process. The most optimal method I discovered",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "a time-consuming trial and error process. The most optimal method I discovered",
  2) | Out-Null

# 6. Text change: reorder "one possible improvement I thought of to lessen the workload
#    of data wrangling was utilising" -> "I thought of one possible improvement to lessen
#    the data wrangling workload - utilising"
$d.Content.Find.Execute(
  "In hindsight, one possible improvement I thought of to lessen the workload of data wrangling was utilising R functions",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "In hindsight, I thought of one possible improvement to lessen the data wrangling workload - utilising R functions",
  2) | Out-Null

# 7a. Text change in the final-paragraph section: drop "even further", move "of the"
$d.Content.Find.Execute(
  "Initially, I wanted to improve the current project even further by scraping the content of the questions and running sentiment analysis",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Initially, I wanted to improve the current project by scraping the questions’ content and running sentiment analysis",
  2) | Out-Null

# 7b. Text change: "states of Members of Parliament." -> "states of the Members of Parliament."
$d.Content.Find.Execute(
  "the cognitive and affective states of Members of Parliament. However, I eventually scraped the idea",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "the cognitive and affective states of the Members of Parliament. However, I eventually scraped the idea",
  2) | Out-Null
